$d = $word.ActiveDocument

# 1. Merge "title" run (with spell-check markers) into the preceding run's text,
#    turning "La primera millora al codi serà afegir un title, " into one run.
$d.Content.Find.Execute("La primera millora al codi serà afegir un title, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "La primera millora al codi serà afegir un title, ", 2)

# 2. Add a new list paragraph after the first list item, before the trailing blank paragraph.
$target = $d.Paragraphs(3)
$range = $target.Range
$range.Collapse(0)
$range.InsertParagraphAfter()
$newPara = $d.Paragraphs(4)
$newPara.Range.Text = "El segon pas es la separació del HTML i el CSS, ja que facilita la gestió i actualització dels estils sense afectar l'estructura HTML, millorant la claredat i l'organització del codi."
$newPara.Style = $target.Style
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($target.Range.ListFormat.ListTemplate)
